# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" before the "总计" (Total) sheet and
#    populate it with the per-fund holding data for 2022-Q1 (same layout
#    as the "2021-Q4" sheet).
# 2. Insert a new summary row at the top of the "总计" sheet's data for
#    the 2022-Q1 quarter, shifting the existing rows down.

$wb = $excel.ActiveWorkbook

$q4Sheet = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# Helper: assign a value that must be stored as TEXT even when it looks
# like a number (e.g. "001410", "140.41"), without leaving a stray
# number-format style behind on the cell.
# ---------------------------------------------------------------------
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Create the "2022-Q1" worksheet right before "总计".
#    NOTE: worksheet variables are position-bound in this runtime, so we
#    must look "总计" up again (by name) right before inserting, and
#    again afterwards whenever we need to touch it - capturing it too
#    early would silently alias the newly inserted sheet instead.
# ---------------------------------------------------------------------
$totalSheetBefore = $wb.Worksheets.Item("总计")
$newWs = $wb.Worksheets.Add($totalSheetBefore)
$newWs.Name = "2022-Q1"

# Copy header (B1:H1) and the index-column (A) formatting from the
# "2021-Q4" sheet so the new sheet matches the existing style (s="2").
$q4Sheet.Range("B1:H1").Copy()
$newWs.Range("B1:H1").PasteSpecial(-4122)

$q4Sheet.Range("A2").Copy()
$newWs.Range("A2:A14").PasteSpecial(-4122)

# Header row
$newWs.Range("B1").Value = "基金代码"
$newWs.Range("C1").Value = "基金名称"
$newWs.Range("D1").Value = "基金规模"
$newWs.Range("E1").Value = "股票总仓位"
$newWs.Range("F1").Value = "仓位占比"
$newWs.Range("G1").Value = "持有市值(亿元)"
$newWs.Range("H1").Value = "仓位排名"

# Data rows
$rows = @(
    @(0,  "001410", "信达澳银新能源产业股票",                 "140.41", "92.06", "1.53", "2.1483", 7),
    @(1,  "012608", "信达澳银领先智选混合型证券投资基金",        "38.78",  "90.57", "1.52", "0.5895", 7),
    @(2,  "011188", "信达澳银星奕混合A",                      "35.73",  "92.14", "1.52", "0.5431", 7),
    @(3,  "006257", "信达澳银先进智造股票",                    "21.53",  "93.84", "1.56", "0.3359", 7),
    @(4,  "007484", "信达澳银核心科技混合",                    "16.02",  "80.50", "1.27", "0.2035", 10),
    @(5,  "009511", "信达澳银研究优选混合",                    "9.41",   "92.12", "1.52", "0.1430", 7),
    @(6,  "011223", "信达澳银星奕混合C",                      "9.21",   "92.14", "1.52", "0.1400", 7),
    @(7,  "009437", "信达澳银科技创新一年定期开放混合A",         "5.57",   "94.43", "1.53", "0.0852", 7),
    @(8,  "003493", "申万菱信安鑫优选混合A",                   "5.34",   "25.97", "0.91", "0.0486", 10),
    @(9,  "004265", "金鹰民丰回报定期开放混合",                 "6.57",   "28.61", "0.70", "0.0460", 8),
    @(10, "003512", "申万菱信安鑫优选混合C",                   "1.82",   "25.97", "0.91", "0.0166", 10),
    @(11, "009438", "信达澳银科技创新一年定期开放混合C",         "0.84",   "94.43", "1.53", "0.0129", 7),
    @(12, "970050", "东海海睿锐意3个月定开灵活配置混合",         "0.17",   "78.42", "2.69", "0.0046", 10)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $newWs.Range("A$r").Value = $row[0]

    Set-TextValue $newWs.Range("B$r") $row[1]
    Set-TextValue $newWs.Range("C$r") $row[2]
    Set-TextValue $newWs.Range("D$r") $row[3]
    Set-TextValue $newWs.Range("E$r") $row[4]
    Set-TextValue $newWs.Range("F$r") $row[5]
    Set-TextValue $newWs.Range("G$r") $row[6]

    $newWs.Range("H$r").Value = $row[7]
}

# ---------------------------------------------------------------------
# 2. Insert the 2022-Q1 summary row into "总计", pushing existing rows
#    down (2021-Q4 becomes row 3, 2021-Q3 becomes row 4).
#    Re-fetch "总计" by name now that the sheet collection has changed.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Re-apply the index-column style (s="2") to the new A2, copied from A3.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2

# Clear any leftover formatting the row-insert may have applied to B:D.
$totalSheet.Range("B2:D2").Style = "Normal"

$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 13
$totalSheet.Range("D2").Value = 4.32
